# The trace report was re-run against the WCS search (new pull completed
# 06/22/2023 11:10:13 EDT instead of 06/21/2023 09:18:52 EDT), which refreshed
# every car-event row in the table. Apply the refreshed data to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds the free-text description / completion timestamp for the search.
$ws.Range('A1').Value = 'Description unknown, completed 06/22/2023 11:10:13 EDT, by WPJTOWN1.The search returned: 8 events.'

# Columns (A:O) for each data row:
# Initial, Number, Location City, State, Month, Day, Time, Event,
# Train ID, Destination City, State, Gross Weight, Tare Weight, Net Weight, Car_no
$rows = @(
    @('MWCX', 100715, 'BIRMINGHAM', 'AL', 6, 21, 1610, 'Arrive In-Transit', ''      , 'LOVELAND',  'CO', 267600, 73600, 194000, 'MWCX100715'),
    @('MWCX', 102553, 'DENVER',     'CO', 6, 21, 1323, 'Arrive In-Transit', 'HKCKDE', 'LOVELAND',  'CO', 281050, 73400, 207650, 'MWCX102553'),
    @('ITFX',   9728, 'JOHNSTOWN',  'CO', 6,  1, 1812, 'Placed Actual',     ''      , 'JOHNSTOWN', 'CO', 202950,     0, 202950, 'ITFX9728'),
    @('MWCX', 102276, 'JOHNSTOWN',  'CO', 6, 12, 1304, 'Placed Actual',     ''      , 'LOVELAND',  'CO', 280350, 78900, 201450, 'MWCX102276'),
    @('MWCX', 102166, 'JOHNSTOWN',  'CO', 6, 12, 1304, 'Placed Actual',     ''      , 'LOVELAND',  'CO', 282400, 82000, 200400, 'MWCX102166'),
    @('MWCX', 102330, 'JOHNSTOWN',  'CO', 6, 15, 1435, 'Placed Actual',     ''      , 'LOVELAND',  'CO', 284850, 79300, 205550, 'MWCX102330'),
    @('MWCX', 102328, 'MEMPHIS',    'TN', 6, 20, 1950, 'Bad Order',         'L 000' , 'LOVELAND',  'CO', 280550, 79500, 201050, 'MWCX102328'),
    @('MWCX', 100705, 'NETTLETON',  'MS', 6, 21, 2000, 'Arrive In-Transit', 'YAMO10', 'LOVELAND',  'CO', 267061, 72200, 194861, 'MWCX100705')
)

$numRows = $rows.Count
$numCols = $rows[0].Count
$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

# The refreshed results now occupy rows 3-10 (same A:O range as before).
$ws.Range('A3:O10').Value = $data
